# "Generate Report for Handback" - CI run that marks the localized files as
# handed back (in sync with en-US) and records the handback target/file/time
# for each locale sheet, widening the columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: was "Ready for handoff", now reflects the handback ---
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$zhcn.Range("C2").Value = $statusText
$dede.Range("C2").Value = $statusText

# Overview columns E/F (zh-cn / de-de status) grow to fit the longer text
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# --- zh-cn sheet: fill in handback target/file/datetime ---
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5799ee0116bc7db7bc3c7b530636ac7cfb6f0dc/e2e/5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.md", [System.Reflection.Missing]::Value, "5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.md", "5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.md")
# match the workbook's existing (non-theme) hyperlink font: underlined cornflower blue
$zhcn.Range("I2").Font.ThemeColor = [System.Reflection.Missing]::Value
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("I2").Font.Underline = $true

$zhcn.Range("J2").Value = "5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.ad206ce02866240825132e03ee6a42447766c14c.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-17 04:53:03"

# --- de-de sheet: fill in handback target/file/datetime ---
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a5799ee0116bc7db7bc3c7b530636ac7cfb6f0dc/e2e/5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.md", [System.Reflection.Missing]::Value, "5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.md", "5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.md")
# match the workbook's existing (non-theme) hyperlink font: underlined cornflower blue
$dede.Range("I2").Font.ThemeColor = [System.Reflection.Missing]::Value
$dede.Range("I2").Font.Color = 15570276
$dede.Range("I2").Font.Underline = $true

$dede.Range("J2").Value = "5e5605a2-dd5c-4ecb-9e5b-1d73fd23d5ee.ad206ce02866240825132e03ee6a42447766c14c.de-de.xlf"
$dede.Range("K2").Value = "2016-08-17 04:53:12"
